$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$cell = $ws.Range("B11")

# Stash the cell's current formatting on an unused scratch row far away
# from the real data so we can restore it after forcing a text (not
# numeric) value into B11.
$scratch = $ws.Range("Z50")
$cell.Copy($scratch)

# Assigning the bare digit string "1" would otherwise be auto-detected
# as a number; switching to a text number format first makes Excel
# store it as literal text "1" (a new shared string), matching the
# change from "R40" to "1" as the cell's displayed text.
$cell.NumberFormat = "@"
$cell.Value = "1"

# Restore the original look (borders/fill/font/alignment/number format)
# of B11 without touching the text value we just set.
$scratch.Copy()
$cell.PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Remove the scratch row entirely so it leaves no trace (no stray
# cell, no change to the sheet's used range/dimension).
$scratch.EntireRow.Delete() | Out-Null
